$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 162.9941356354778
$ws.Range("R2").Value = 504.5149151680536
$ws.Range("S2").Value = 2807.261316242745

$ws.Range("M3").Value = 0
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 228.7895586469373
$ws.Range("R3").Value = 2461.222272034012
$ws.Range("S3").Value = 2720.727186006326

$ws.Range("M4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 219.4935176313175
$ws.Range("R4").Value = 5403.955159823739
$ws.Range("S4").Value = 2727.070397138864

$ws.Range("A5").Value = 3
$ws.Range("M5").Value = 0
$ws.Range("P5").Value = 0
$ws.Range("Q5").Value = 219.4935176313175
$ws.Range("R5").Value = 5403.955159823739
$ws.Range("S5").Value = 2727.070397138864

$ws.Range("A6").Value = 4
$ws.Range("M6").Value = 0
$ws.Range("P6").Value = 0
$ws.Range("Q6").Value = 219.4935176313175
$ws.Range("R6").Value = 5403.955159823739
$ws.Range("S6").Value = 2727.070397138864

$ws.Range("A7").Value = 5
$ws.Range("M7").Value = 0
$ws.Range("P7").Value = 0
$ws.Range("Q7").Value = 219.4935176313175
$ws.Range("R7").Value = 5403.955159823739
$ws.Range("S7").Value = 2727.070397138864

$ws.Range("A8").Value = 6
$ws.Range("M8").Value = 0
$ws.Range("P8").Value = 0
$ws.Range("Q8").Value = 219.4935176313175
$ws.Range("R8").Value = 5403.955159823739
$ws.Range("S8").Value = 2727.070397138864

$ws.Range("A9").Value = 7
$ws.Range("M9").Value = 0
$ws.Range("P9").Value = 0
$ws.Range("Q9").Value = 219.4935176313175
$ws.Range("R9").Value = 5403.955159823739
$ws.Range("S9").Value = 2727.070397138864

$ws.Range("A10").Value = 8
$ws.Range("M10").Value = 0
$ws.Range("P10").Value = 0
$ws.Range("Q10").Value = 219.4935176313175
$ws.Range("R10").Value = 5403.955159823739
$ws.Range("S10").Value = 2727.070397138864

$ws.Range("A11").Value = 9
$ws.Range("M11").Value = 0
$ws.Range("P11").Value = 0
$ws.Range("Q11").Value = 219.4935176313175
$ws.Range("R11").Value = 5403.955159823739
$ws.Range("S11").Value = 2727.070397138864

$ws.Range("A12").Value = 10
$ws.Range("M12").Value = 0
$ws.Range("P12").Value = 0
$ws.Range("Q12").Value = 219.4935176313175
$ws.Range("R12").Value = 5403.955159823739
$ws.Range("S12").Value = 2727.070397138864

$ws.Range("A13").Value = 11
$ws.Range("M13").Value = 0
$ws.Range("P13").Value = 0
$ws.Range("Q13").Value = 219.4935176313175
$ws.Range("R13").Value = 5403.955159823739
$ws.Range("S13").Value = 2727.070397138864

$ws.Range("A14").Value = 12
$ws.Range("M14").Value = 0
$ws.Range("P14").Value = 0
$ws.Range("Q14").Value = 219.4935176313175
$ws.Range("R14").Value = 5403.955159823739
$ws.Range("S14").Value = 2727.070397138864

$ws.Range("A15").Value = 13
$ws.Range("M15").Value = 0
$ws.Range("P15").Value = 0
$ws.Range("Q15").Value = 219.4935176313175
$ws.Range("R15").Value = 5403.955159823739
$ws.Range("S15").Value = 2727.070397138864

$ws.Range("A16").Value = 14
$ws.Range("M16").Value = 0
$ws.Range("P16").Value = 0
$ws.Range("Q16").Value = 219.4935176313175
$ws.Range("R16").Value = 5403.955159823739
$ws.Range("S16").Value = 2727.070397138864

$ws.Range("A17").Value = 15
$ws.Range("M17").Value = 0
$ws.Range("P17").Value = 0
$ws.Range("Q17").Value = 219.4935176313175
$ws.Range("R17").Value = 5403.955159823739
$ws.Range("S17").Value = 2727.070397138864

$ws.Range("A18").Value = 16
$ws.Range("M18").Value = 0
$ws.Range("P18").Value = 0
$ws.Range("Q18").Value = 219.4935176313175
$ws.Range("R18").Value = 5403.955159823739
$ws.Range("S18").Value = 2727.070397138864

$ws.Range("A19").Value = 17
$ws.Range("M19").Value = 0
$ws.Range("P19").Value = 0
$ws.Range("Q19").Value = 219.4935176313175
$ws.Range("R19").Value = 5403.955159823739
$ws.Range("S19").Value = 2727.070397138864

$ws.Range("A20").Value = 18
$ws.Range("M20").Value = 0
$ws.Range("P20").Value = 0
$ws.Range("Q20").Value = 219.4935176313175
$ws.Range("R20").Value = 5403.955159823739
$ws.Range("S20").Value = 2727.070397138864

$ws.Range("A21").Value = 19
$ws.Range("M21").Value = 0
$ws.Range("P21").Value = 0
$ws.Range("Q21").Value = 219.4935176313175
$ws.Range("R21").Value = 5403.955159823739
$ws.Range("S21").Value = 2727.070397138864

$ws.Range("A22").Value = 20
$ws.Range("M22").Value = 0
$ws.Range("P22").Value = 0
$ws.Range("Q22").Value = 219.4935176313175
$ws.Range("R22").Value = 5403.955159823739
$ws.Range("S22").Value = 2727.070397138864

$ws.Range("A23").Value = 21
$ws.Range("M23").Value = 0
$ws.Range("P23").Value = 0
$ws.Range("Q23").Value = 219.4935176313175
$ws.Range("R23").Value = 5403.955159823739
$ws.Range("S23").Value = 2727.070397138864

$ws.Range("A24").Value = 22
$ws.Range("M24").Value = 0
$ws.Range("P24").Value = 0
$ws.Range("Q24").Value = 219.4935176313175
$ws.Range("R24").Value = 5403.955159823739
$ws.Range("S24").Value = 2727.070397138864

$ws.Range("A25").Value = 23
$ws.Range("M25").Value = 0
$ws.Range("P25").Value = 0
$ws.Range("Q25").Value = 219.4935176313175
$ws.Range("R25").Value = 5403.955159823739
$ws.Range("S25").Value = 2727.070397138864

$ws.Range("A26").Value = 24
$ws.Range("M26").Value = 0
$ws.Range("P26").Value = 0
$ws.Range("Q26").Value = 219.4935176313175
$ws.Range("R26").Value = 5403.955159823739
$ws.Range("S26").Value = 2727.070397138864
